# "day 6 intorduction to bootstrap"
# The active workbook's active sheet is "Week 2" (tabSelected=1 / activeTab index 1).
# This day's edits fill in the "Task 2" (column E) attendance marks ("p") for most
# students, set the Task 2 due-date header (E3), and move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make sure we're working on the sheet that was active when the file was saved.
$ws.Activate()

# Header row: Task 2 date (one day after the Task 1 date already in D3).
$ws.Range("E3").Value = 45454

# Column E ("Task 2") gets the same "p" mark already present in column D for
# every student row except the handful that stayed blank (rows 7, 11, 15, 16,
# 19, 20).
$rowsToMark = @(4, 5, 6, 8, 9, 10, 12, 13, 14, 17, 18, 21, 22, 23, 24, 25)
foreach ($r in $rowsToMark) {
    $ws.Cells.Item($r, 5).Value = "p"
}

# Move the active selection / scrolled view to match where the author ended up.
$ws.Range("I14").Select()
$excel.ActiveWindow.ScrollRow = 6
$excel.ActiveWindow.ScrollColumn = 1
